# AlexNet architecture parameter correction
# Corrects conv1..conv5 layer dimensions and inserts the missing conv5 layer
# (plus the two accidentally-missing dense2/dense-final/softmax rows at the
# tail) so the MultAddsComp sheet reflects the real AlexNet architecture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room in the table: one new row for a missing "max pooling"
#    layer after conv2 (becomes row 8), and three new rows for the
#    missing dense2 / relu / dense final / softmax tail entries
#    (rows 19-21) before the "total" row.
# ---------------------------------------------------------------------
$ws.Rows("8:8").Insert()
$ws.Rows("19:21").Insert()

# ---------------------------------------------------------------------
# 2) conv1 block (rows 3-4): 48 -> 96 feature maps, 5x5 -> 11x11 kernel
# ---------------------------------------------------------------------
$ws.Range("D3").Value = 96
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = 11
$ws.Range("D4").Value = 96

# ---------------------------------------------------------------------
# 3) max pooling after conv1 (row 5): 13x13x192 -> 27x27x96
# ---------------------------------------------------------------------
$ws.Range("B5").Value = 27
$ws.Range("C5").Value = 27
$ws.Range("D5").Value = 96

# ---------------------------------------------------------------------
# 4) conv2 (row 6) + relu (row 7): 13x13x192 -> 27x27x256
# ---------------------------------------------------------------------
$ws.Range("B6").Value = 27
$ws.Range("C6").Value = 27
$ws.Range("D6").Value = 256
$ws.Range("B7").Value = 27
$ws.Range("C7").Value = 27
$ws.Range("D7").Value = 256

# ---------------------------------------------------------------------
# 5) row 8 (new): max pooling after conv2 -> 13x13x256
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "max pooling"
$ws.Range("B8").Value = 13
$ws.Range("C8").Value = 13
$ws.Range("D8").Value = 256
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 3

# ---------------------------------------------------------------------
# 6) row 9: conv3, 13x13x192 -> 13x13x384; clear the stale G formula
#    (conv3 no longer reports MultAdds in its own row G column)
# ---------------------------------------------------------------------
$ws.Range("D9").Value = 384
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 3
$ws.Range("G9").Clear()

# ---------------------------------------------------------------------
# 7) row 10: relu, 13x13x192 -> 13x13x384; add a simple param-count
#    formula in G (no MultAdds for relu)
# ---------------------------------------------------------------------
$ws.Range("D10").Value = 384
$ws.Range("G10").Formula = "=B10*C10*D10"

# ---------------------------------------------------------------------
# 8) row 11: conv4, 13x13x128 -> 13x13x384; clear stale G formula
# ---------------------------------------------------------------------
$ws.Range("D11").Value = 384
$ws.Range("G11").Clear()

# ---------------------------------------------------------------------
# 9) row 12: relu, 13x13x192 -> 13x13x384 (label already "relu")
# ---------------------------------------------------------------------
$ws.Range("D12").Value = 384

# ---------------------------------------------------------------------
# 10) row 13: new conv5 layer, 13x13x384 -> 13x13x256
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "conv5"
$ws.Range("B13").Value = 13
$ws.Range("C13").Value = 13
$ws.Range("D13").Value = 256
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 3
$ws.Range("H13").Clear()
$ws.Range("I13").Clear()
$ws.Range("J13").Clear()

# ---------------------------------------------------------------------
# 11) row 14: relu, 1x1x4096 -> 13x13x256
# ---------------------------------------------------------------------
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = 13
$ws.Range("D14").Value = 256

# ---------------------------------------------------------------------
# 12) row 15: max pooling (after conv5), 1x1x4096 -> 4x4x256
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "max pooling"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 256
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 3
$ws.Range("G15").Clear()
$ws.Range("H15").Clear()
$ws.Range("I15").Clear()
$ws.Range("J15").Clear()

# ---------------------------------------------------------------------
# 13) row 16: dense1, E/F 1 -> 13 (now the real fully-connected layer)
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "dense1"
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 13
$ws.Range("G16").Formula = "=2*E16*F16*D12*D16*B16*C16"
$ws.Range("H16").Value = 150528
$ws.Range("H16").Font.Color = 0
$ws.Range("I16").Formula = "=E16*F16*D12*D16"
$ws.Range("J16").Formula = "=D16*C16*B16"

# ---------------------------------------------------------------------
# 14) row 17: relu after dense1 (static 1x1x1000 shape kept)
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "relu"
$ws.Range("G17").Clear()
$ws.Range("H17").Clear()
$ws.Range("I17").Clear()
$ws.Range("J17").Clear()

# ---------------------------------------------------------------------
# 15) row 18: dense2, 1x1x1000 -> 1x1x4096
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "dense2"
$ws.Range("D18").Value = 4096
$ws.Range("G18").Value = 21780000
$ws.Range("G18").Font.Color = 0
$ws.Range("H18").Value = 150528
$ws.Range("H18").Font.Color = 0
$ws.Range("I18").Formula = "=E18*F18*D17*D18"
$ws.Range("J18").Formula = "=D18*C18*B18"

# ---------------------------------------------------------------------
# 16) row 19 (new): relu, 1x1x4096
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "relu"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 4096
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 1

# ---------------------------------------------------------------------
# 17) row 20 (new): dense final, 1x1x1000
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "dense final"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1000
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 21780000
$ws.Range("G20").Font.Color = 0
$ws.Range("H20").Value = 150528
$ws.Range("H20").Font.Color = 0
$ws.Range("I20").Formula = "=E20*F20*D19*D20"
$ws.Range("J20").Formula = "=D20*C20*B20"

# ---------------------------------------------------------------------
# 18) row 21 (new): softmax, 1x1x1000
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "softmax"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1000
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1

# ---------------------------------------------------------------------
# 19) row 22: total row - extend SUM ranges to cover the new rows
# ---------------------------------------------------------------------
$ws.Range("G22").Formula = "=SUM(G2:G21)"
$ws.Range("H22").Formula = "=SUM(H2:H21)"
$ws.Range("I22").Formula = "=SUM(I2:I21)"
$ws.Range("J22").Formula = "=SUM(J2:J21)"

Write-Host "AlexNet complexity table corrected"
